$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the header cell / table column from "Nom Comarca" to "Nom"
$ws.Range("B3").Value = "Nom"
